$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 93.6022758618435
$ws.Range("F2").Value = 14738
$ws.Range("G2").Value = 157.453436514228
$ws.Range("I2").Value = 0.091
$ws.Range("J2").Value = 7333
$ws.Range("K2").Value = 5746
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2364
$ws.Range("O2").Value = 51.4409
$ws.Range("P2").Value = 8.05941363323641
$ws.Range("Q2").Value = 14.6877224931919
$ws.Range("R2").Value = 51
$ws.Range("S2").Value = 0.54
$ws.Range("E3").Value = 96.6011004628816
$ws.Range("F3").Value = 2871
$ws.Range("G3").Value = 29.7201583236949
$ws.Range("J3").Value = 2626
$ws.Range("K3").Value = 21
$ws.Range("L3").Value = 216
$ws.Range("N3").Value = 76
$ws.Range("O3").Value = 89.3893333333333
$ws.Range("P3").Value = 2.97911209099006
$ws.Range("Q3").Value = 5.20403060332032
$ws.Range("R3").Value = "#NUM!"
$ws.Range("S3").Value = "#NUM!"
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 74.4371247489178
$ws.Range("F4").Value = 25972
$ws.Range("G4").Value = 348.91191844937
$ws.Range("I4").Value = 0.45
$ws.Range("J4").Value = 17965
$ws.Range("K4").Value = 4775
$ws.Range("L4").Value = 63
$ws.Range("M4").Value = 207
$ws.Range("N4").Value = 4397
$ws.Range("O4").Value = 66.583
$ws.Range("P4").Value = 4.70032783052408
$ws.Range("Q4").Value = 10.0015149709033
$ws.Range("R4").Value = 35.7894736842105
$ws.Range("D5").Value = 76
$ws.Range("E5").Value = 55.4265564192303
$ws.Range("F5").Value = 76827
$ws.Range("G5").Value = 1386.10451313091
$ws.Range("I5").Value = 0.11
$ws.Range("J5").Value = 21381
$ws.Range("K5").Value = 51088
$ws.Range("L5").Value = 162
$ws.Range("M5").Value = 876
$ws.Range("N5").Value = 2960
$ws.Range("O5").Value = 42.7021746031746
$ws.Range("P5").Value = 12.1891076133896
$ws.Range("Q5").Value = 14.0619655080493
$ws.Range("S5").Value = 0.37
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 49.0416132262602
$ws.Range("F6").Value = 71590
$ws.Range("G6").Value = 1459.78069011942
$ws.Range("I6").Value = 0.25
$ws.Range("J6").Value = 66753
$ws.Range("K6").Value = 1838
$ws.Range("L6").Value = 44
$ws.Range("M6").Value = 778
$ws.Range("N6").Value = 1298
$ws.Range("O6").Value = 84.0847105263158
$ws.Range("P6").Value = 4.02069627852412
$ws.Range("Q6").Value = 2.70049621750537
$ws.Range("R6").Value = 29.375
$ws.Range("S6").Value = 0.365625
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 155.182033671545
$ws.Range("F7").Value = 9011
$ws.Range("G7").Value = 58.067289020535
$ws.Range("J7").Value = 8759
$ws.Range("K7").Value = 24
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 28
$ws.Range("N7").Value = 215
$ws.Range("O7").Value = 62.2262857142857
$ws.Range("P7").Value = 11.0061003514899
$ws.Range("Q7").Value = 5.11730037213777
$ws.Range("R7").Value = "#NUM!"
$ws.Range("S7").Value = "#NUM!"
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 150.568236356283
$ws.Range("F8").Value = 17483
$ws.Range("G8").Value = 116.113467375886
$ws.Range("J8").Value = 9871
$ws.Range("K8").Value = 6032
$ws.Range("L8").Value = 24
$ws.Range("M8").Value = 110
$ws.Range("N8").Value = 1262
$ws.Range("O8").Value = 71.8068333333333
$ws.Range("P8").Value = 6.16883736293026
$ws.Range("Q8").Value = 6.59066162116396
$ws.Range("R8").Value = 58
$ws.Range("S8").Value = 0.46
$ws.Range("E9").Value = 179.105185543933
$ws.Range("F9").Value = 20205
$ws.Range("G9").Value = 112.810804101726
$ws.Range("J9").Value = 8866
$ws.Range("K9").Value = 9730
$ws.Range("L9").Value = 24
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 1310
$ws.Range("O9").Value = 65.7231111111111
$ws.Range("P9").Value = 8.3655822976522
$ws.Range("Q9").Value = 5.49319235766674
$ws.Range("R9").Value = 60
$ws.Range("S9").Value = 0.5
$ws.Range("D10").Value = 32
$ws.Range("E10").Value = 41.167605010108
$ws.Range("F10").Value = 33868
$ws.Range("G10").Value = 822.685701334441
$ws.Range("I10").Value = 0.062
$ws.Range("J10").Value = 29237
$ws.Range("K10").Value = 1493
$ws.Range("L10").Value = 142
$ws.Range("M10").Value = 862
$ws.Range("N10").Value = 1384
$ws.Range("O10").Value = 67.3829333333333
$ws.Range("P10").Value = 6.85952029905608
$ws.Range("Q10").Value = 8.30876749341512
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 315.918408936201
$ws.Range("F11").Value = 18434
$ws.Range("G11").Value = 58.3505091142781
$ws.Range("I11").Value = 0.91
$ws.Range("J11").Value = 16302
$ws.Range("K11").Value = 424
$ws.Range("L11").Value = 9
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 8255
$ws.Range("O11").Value = 85.7849
$ws.Range("P11").Value = 3.69977641559955
$ws.Range("Q11").Value = 3.38793267483305
$ws.Range("R11").Value = 20
